$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings / error messages for the newly added rows
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Argumentos de ruta invalidos"
$ws.Range("C6").Value = 400

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Ruta de api invalida"
$ws.Range("C7").Value = 404

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "server internal error"
$ws.Range("C8").Value = 500

# New column C: HTTP status codes for the pre-existing rows
$ws.Range("C1").Value = 200
$ws.Range("C3").Value = 400
$ws.Range("C4").Value = 404

# Update the active selection to match the saved state (C7)
$ws.Range("C7").Select()
